$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.234.50"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "3.319.14"
$ws.Range("E3").Value = "  -3.10%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.32"
$ws.Range("E5").Value = "  -2.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "650.47"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  -7.49%  "

$ws.Range("E8").Value = "  -2.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.984"
$ws.Range("E10").Value = "  -7.53%  "

$ws.Range("D11").Value = "3.316.19"
$ws.Range("E11").Value = "  -3.00%  "

$ws.Range("E12").Value = "  -3.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.97"
$ws.Range("E13").Value = "  -4.97%  "

$ws.Range("D14").Value = "96.001.87"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.06"
$ws.Range("E15").Value = "  -4.56%  "

$ws.Range("E16").Value = "  -4.46%  "

$ws.Range("D17").Value = "3.936.59"
$ws.Range("E17").Value = "  -3.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.46"
$ws.Range("E18").Value = "  -2.70%  "

$ws.Range("D19").Value = "3.316.97"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.534"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.98"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "501.68"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.44"
$ws.Range("E23").Value = "  -4.73%  "

$ws.Range("E24").Value = "  -3.56%  "

$ws.Range("E25").Value = "  -4.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("E26").Value = "  +7.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.62"
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.98"
$ws.Range("E28").Value = "  -6.55%  "

$ws.Range("E29").Value = "  -9.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.95"
$ws.Range("E31").Value = "  -4.32%  "

$ws.Range("E32").Value = "  -7.15%  "

$ws.Range("E33").Value = "  +7.54%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  -6.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.86"
$ws.Range("E36").Value = "  -7.11%  "

$ws.Range("E37").Value = "  +2.09%  "

$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("E40").Value = "  -3.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "503.01"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("E43").Value = "  +1.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.824"
$ws.Range("E44").Value = "  -4.56%  "

$ws.Range("E45").Value = "  -1.58%  "

$ws.Range("E46").Value = "  +4.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.46"
$ws.Range("E47").Value = "  -0.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.30"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.93"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.11"
$ws.Range("E50").Value = "  -5.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "161.55"
$ws.Range("E51").Value = "  -0.47%  "
